# Pinout.xlsx - "Updated PinAssignment on SAM4S-Board"
#
# Adds a new USART bus column (K8/D11/G12), marks the new W25Q128.SPI
# flash connections (R11/O12/O14/R14), and highlights the pins that are
# now in use (blue fill) across the pinout table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blue   = 12611584   # RGB(0,112,192)  -> FF0070C0
$orange = 49407       # RGB(255,192,0)  -> FFFFC000 (existing "bus" highlight colour)
$xlCenter = -4108     # xlCenter
$xlEdgeRight = 10     # xlEdgeRight
$xlContinuous = 1     # xlContinuous
$xlNone = 0          # xlLineStyleNone (handy for "clear borders, then re-add one edge")

# ---------------------------------------------------------------
# Row 8: new "USART" bus column header cell (K8), styled like the
# other bus-label cells (K5/K6/K7) but with the new blue colour and
# a right-hand border only.
# ---------------------------------------------------------------
$ws.Range("K8").Value = "USART"
$ws.Range("K8").Interior.Color = $blue
$ws.Range("K8").HorizontalAlignment = $xlCenter
$ws.Range("K8").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

# Highlight the pins newly wired to USART
$ws.Range("I8").Interior.Color = $blue
$ws.Range("L8").Interior.Color = $blue

# ---------------------------------------------------------------
# Row 9
# ---------------------------------------------------------------
$ws.Range("L9").Interior.Color = $blue
$ws.Range("S9").Interior.Color = $blue

# ---------------------------------------------------------------
# Row 10
# ---------------------------------------------------------------
$ws.Range("A10").Interior.Color = $blue
$ws.Range("H10").Interior.Color = $blue

# ---------------------------------------------------------------
# Row 11
# ---------------------------------------------------------------
$ws.Range("A11").Interior.Color = $blue

# VESC.USART label (blue, centered, right border only)
$ws.Range("D11").Borders.LineStyle = $xlNone
$ws.Range("D11").Value = "VESC.USART"
$ws.Range("D11").Interior.Color = $blue
$ws.Range("D11").HorizontalAlignment = $xlCenter
$ws.Range("D11").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

# W25Q128.SPI flash wiring (reuses the existing orange "bus" fill)
$ws.Range("R11").Value = "W25Q128.SPI"
$ws.Range("R11").Interior.Color = $orange

# ---------------------------------------------------------------
# Row 12
# ---------------------------------------------------------------
$ws.Range("G12").Borders.LineStyle = $xlNone
$ws.Range("G12").Value = "VESC.USART"
$ws.Range("G12").Interior.Color = $blue
$ws.Range("G12").HorizontalAlignment = $xlCenter
$ws.Range("G12").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

$ws.Range("H12").Interior.Color = $blue

$ws.Range("O12").Value = "W25Q128.SPI"
$ws.Range("O12").Interior.Color = $orange

$ws.Range("S12").Interior.Color = $blue

# ---------------------------------------------------------------
# Row 13
# ---------------------------------------------------------------
$ws.Range("A13").Interior.Color = $blue
$ws.Range("H13").Interior.Color = $blue
$ws.Range("L13").Interior.Color = $blue
$ws.Range("S13").Interior.Color = $blue

# ---------------------------------------------------------------
# Row 14
# ---------------------------------------------------------------
$ws.Range("O14").Value = "W25Q128.SPI"
$ws.Range("O14").Interior.Color = $orange

$ws.Range("R14").Value = "W25Q128.SPI"
$ws.Range("R14").Interior.Color = $orange

# ---------------------------------------------------------------
# Final selection, as left by the author after the edit.
# ---------------------------------------------------------------
$ws.Range("L21").Select()
